$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the species-record data held in rows 7-10 (columns A,B,D,E,F,G,H,Q,R)
# while leaving the other, unchanged columns (C, I, J, K, N, P, S, T, U, V, W, Y, Z, AA, AB,
# AD, AE, AF, AG, AT, AW, AX, AY) untouched:
#   new row 7  <- old row 10 data
#   new row 8  <- old row 9 data
#   new row 9  <- old row 7 data
#   new row 10 <- old row 8 data

# Row 7 (becomes old row 10's record)
$ws.Range("A7").Value  = 111670690
$ws.Range("B7").Value  = 96348
$ws.Range("D7").Value  = "VU"
$ws.Range("E7").Value  = 220787
$ws.Range("F7").Value  = "Knärot"
$ws.Range("G7").Value  = "Goodyera repens"
$ws.Range("H7").Value  = "(L.) R. Br."
$ws.Range("L7").Value  = ""
$ws.Range("Q7").Value  = 557809.1117697239
$ws.Range("R7").Value  = 7067699.199123298

# Row 8 (becomes old row 9's record)
$ws.Range("A8").Value  = 111671165
$ws.Range("B8").Value  = 78578
$ws.Range("D8").Value  = "NT"
$ws.Range("E8").Value  = 6458
$ws.Range("F8").Value  = "Lunglav"
$ws.Range("G8").Value  = "Lobaria pulmonaria"
$ws.Range("H8").Value  = "(L.) Hoffm."
$ws.Range("Q8").Value  = 558014.2710882163
$ws.Range("R8").Value  = 7067448.175823289

# Row 9 (becomes old row 7's record)
$ws.Range("A9").Value  = 111670912
$ws.Range("B9").Value  = 78578
$ws.Range("D9").Value  = "NT"
$ws.Range("E9").Value  = 6458
$ws.Range("F9").Value  = "Lunglav"
$ws.Range("G9").Value  = "Lobaria pulmonaria"
$ws.Range("H9").Value  = "(L.) Hoffm."
$ws.Range("Q9").Value  = 557803.3534448177
$ws.Range("R9").Value  = 7067771.317107533

# Row 10 (becomes old row 8's record)
$ws.Range("A10").Value = 111671159
$ws.Range("B10").Value = 81248
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 1312
$ws.Range("F10").Value = "Gammelgransskål"
$ws.Range("G10").Value = "Pseudographis pinicola"
$ws.Range("H10").Value = "(Nyl.) Rehm"
$ws.Range("L10").Value = ""
$ws.Range("Q10").Value = 558006.0394731871
$ws.Range("R10").Value = 7067389.087574247
